$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Source Type: Statistical Institution" Enterprises density row (row 11)
# Micro, SMEs, MSMEs density values get slightly more precise figures.
# Values are stored as text (not numbers) in the workbook, so we write them via
# a text formula and Paste Special (values only) to avoid Excel auto-converting
# the numeric-looking strings into number cells.

function Set-TextValue {
    param($range, [string]$text)
    $range.Formula = '="' + $text + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

Set-TextValue $ws.Range("B11") "26.71"
Set-TextValue $ws.Range("C11") "1.47"
Set-TextValue $ws.Range("D11") "28.17"

# Update the "Source Type: SME Associations" Enterprises density row (row 28)
Set-TextValue $ws.Range("B28") "21.27"
Set-TextValue $ws.Range("C28") "8.15"
Set-TextValue $ws.Range("D28") "29.42"
